# Refresh the "cryptos" price/volume snapshot (Price column D, Volume(1h)
# column E) for rows 2-51 to the values captured by the latest GitHub
# Actions run. Numeric-looking Price strings get NumberFormat "@" (Text)
# applied first so Excel stores them verbatim (e.g. "1.000", "0.9999")
# instead of silently coercing them to real numbers, matching the
# original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.310.01'
$ws.Range('E2').Value = '  +0.09%  '

$ws.Range('D3').Value = '1.930.44'
$ws.Range('E3').Value = '  +0.00%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7559'
$ws.Range('E5').Value = '  +5.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.78'
$ws.Range('E6').Value = '  -2.23%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').Value = '  -0.18%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.12'
$ws.Range('E8').Value = '  +2.73%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3185'
$ws.Range('E9').Value = '  -0.93%  '

$ws.Range('E10').Value = '  -1.08%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7793'
$ws.Range('E11').Value = '  -1.62%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08026'
$ws.Range('E12').Value = '  -0.08%  '

$ws.Range('D13').Value = '1.938.24'
$ws.Range('E13').Value = '  +0.58%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.391'
$ws.Range('E14').Value = '  +0.35%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.47'
$ws.Range('E15').Value = '  -1.33%  '

$ws.Range('E16').Value = '  -0.85%  '

$ws.Range('D17').Value = '30.290.62'
$ws.Range('E17').Value = '  +0.05%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '253.49'
$ws.Range('E18').Value = '  -1.42%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.989'
$ws.Range('E19').Value = '  +4.29%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007962'
$ws.Range('E20').Value = '  -1.30%  '

$ws.Range('D21').Value = '2.186.37'
$ws.Range('E21').Value = '  +0.22%  '

$ws.Range('E22').Value = '  -0.02%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.13%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.717'
$ws.Range('E24').Value = '  -1.37%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.516'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.32'
$ws.Range('E26').Value = '  -0.80%  '

$ws.Range('E27').Value = '  -0.52%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1330'
$ws.Range('E28').Value = '  +4.35%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.214'
$ws.Range('E29').Value = '  -2.94%  '

$ws.Range('E30').Value = '  +1.02%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.516'
$ws.Range('E31').Value = '  -0.93%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.413'
$ws.Range('E32').Value = '  +0.33%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.147'
$ws.Range('E33').Value = '  +0.30%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05225'
$ws.Range('E34').Value = '  +1.36%  '

$ws.Range('E35').Value = '  +4.95%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7547'
$ws.Range('E36').Value = '  +1.55%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.791'
$ws.Range('E37').Value = '  +0.89%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01951'
$ws.Range('E38').Value = '  -0.19%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.808'
$ws.Range('E39').Value = '  +0.00%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.45'
$ws.Range('E40').Value = '  +0.99%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.506'
$ws.Range('E41').Value = '  +2.29%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4491'
$ws.Range('E42').Value = '  -0.07%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.970'
$ws.Range('E43').Value = '  -0.92%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9999'
$ws.Range('E44').Value = '  -0.14%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8366'
$ws.Range('E45').Value = '  -0.65%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.04'
$ws.Range('E46').Value = '  +3.35%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.21'
$ws.Range('E47').Value = '  +0.19%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.598'
$ws.Range('E48').Value = '  +2.19%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.72'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '979.57'
$ws.Range('E50').Value = '  +7.50%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1212'
$ws.Range('E51').Value = '  +6.89%  '

